$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(194, 1).Value = 193
$ws.Cells.Item(194, 2).Value = 1
$ws.Cells.Item(194, 3).Value = "2024-06-19 02:01:45"
$ws.Cells.Item(194, 4).Value = 200
$ws.Cells.Item(194, 5).Value = 18

$ws.Cells.Item(195, 1).Value = 194
$ws.Cells.Item(195, 2).Value = 2
$ws.Cells.Item(195, 3).Value = "2024-06-19 02:01:46"
$ws.Cells.Item(195, 4).Value = 200
$ws.Cells.Item(195, 5).Value = 3
